# Bug fix + add data and menu for full date range
# Fills in the missing Date (A) and MealType (B) columns for rows
# 114-121 (Dec 29 & Dec 30, 2022) and appends four brand new rows
# 122-125 for Dec 31, 2022 (date + meal type only, no coupon data yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells whose formatting we want to reuse:
#  - A2 carries the "date" style (Arial font, date number format)
#  - B2 carries the "MealType" style (Arial font, text)
$dateTemplate = $ws.Range("A2")
$mealTemplate = $ws.Range("B2")

$mealNames = @("BreakFast", "Lunch", "EveningSnacks", "Dinner")

# Dates (as Excel serials) for each block of 4 rows, starting at row 114.
$dates = @(44924, 44924, 44924, 44924, 44925, 44925, 44925, 44925, 44926, 44926, 44926, 44926)

$startRow = 114
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    $aCell = $ws.Cells.Item($row, 1)
    $dateTemplate.Copy($aCell)
    $aCell.NumberFormat = "m/d/yyyy"
    $aCell.Value = $dates[$i]

    $bCell = $ws.Cells.Item($row, 2)
    $mealTemplate.Copy($bCell)
    $bCell.Value = $mealNames[$i % 4]
}
